# The Input_Value sheet held test credentials (URL / UserName / Password)
# in row 2, columns L:N. This edit wipes that stored data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

$ws.Range("L2:N2").ClearContents()

# Leave the selection reflecting the cells that were just cleared.
$ws.Range("L2:N2").Select()

# The header row no longer carries an explicit custom height in the
# saved workbook; let Excel recompute it from the (new) default.
$ws.Rows.Item(1).AutoFit()
